# LeaveBalance 1-10 policies
# - Flip "RunMode" (column C) from yes -> No for rows 12-42 and 84-124
#   on the LeaveBalance sheet.
# - Add a new "dummySheet" after LeaveBalance, seeded with the header
#   row and the first data row (LeaveBalance_01 scenario), mirroring
#   row heights/styles from the source sheet.
# - Update the selection/active-sheet state to match the new layout.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LeaveBalance")

# --- Update RunMode (column C) for rows 12-42 and 84-124: "yes" -> "No" ---
for ($r = 12; $r -le 42; $r++) {
    $ws1.Cells.Item($r, 3).Value = "No"
}
for ($r = 84; $r -le 124; $r++) {
    $ws1.Cells.Item($r, 3).Value = "No"
}

# --- Add the new "dummySheet" worksheet after LeaveBalance ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "dummySheet"

# Copy header (row 1) + first scenario row (row 2) from LeaveBalance,
# preserving cell styles/formatting.
$ws1.Range("A1:V2").Copy($newSheet.Range("A1"))

# Match row heights to the source sheet's wrapped-text rows.
$newSheet.Rows.Item(1).RowHeight = $ws1.Rows.Item(1).RowHeight
$newSheet.Rows.Item(2).RowHeight = $ws1.Rows.Item(2).RowHeight

# --- Selection / view state ---
# LeaveBalance keeps its frozen header pane, but selection now spans
# the full width of the top two rows instead of just column C.
$ws1.Range("A1:XFD2").Select()

# dummySheet becomes the active sheet/tab, selected over the same
# full-width top-two-row range.
$newSheet.Range("A1:XFD2").Select()
$newSheet.Activate()
